# Minor modification of southeast_asia data
# Fill in the diagonal (self-distance = 0) of the distance matrix and
# round the off-diagonal distance values to the nearest integer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("distance")

# Row 3 - Cambodia
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 655
$ws.Range("D3").Value = 1293
$ws.Range("E3").Value = 499
$ws.Range("F3").Value = 463

# Row 4 - Laos
$ws.Range("B4").Value = 655
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 813
$ws.Range("E4").Value = 477
$ws.Range("F4").Value = 339

# Row 5 - Myanmar
$ws.Range("B5").Value = 1293
$ws.Range("C5").Value = 813
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 820
$ws.Range("F5").Value = 1144

# Row 6 - Thailand
$ws.Range("B6").Value = 499
$ws.Range("C6").Value = 477
$ws.Range("D6").Value = 820
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 591

# Row 7 - Vietnam
$ws.Range("B7").Value = 463
$ws.Range("C7").Value = 339
$ws.Range("D7").Value = 1144
$ws.Range("E7").Value = 591
$ws.Range("F7").Value = 0

# Move selection to reflect the author's last cursor position
$ws.Range("F15").Select()
